$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999549410293126,
    0.9990639202211722,
    0.9999580877459376,
    0.9999455938856813,
    0.9999500819956132,
    0.00004206057127643839,
    0.0008737893843822983,
    0.00003033244237326141,
    0.00007925478362657401,
    0.00005479361299991771,
    0.0004236862227404525,
    0.006485412190172526,
    1.0000514959665,
    0.00676150961376217,
    110.1527996107309,
    165.0022117297999
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
